$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# Helper: locate a unique marker substring anywhere in the document and
# replace it with the given text, then nudge character formatting (an
# idempotent Bold toggle) so the replaced span becomes its own <w:r> run
# instead of being silently re-merged into its neighbour. The toggle is a
# no-op on the final formatting (Bold ends up False again) but it forces
# the host to materialise a run boundary, and the run correctly inherits
# the surrounding (already-correct) character formatting.
# ---------------------------------------------------------------------------
function Split-Marker($marker, $replacement) {
    $sr = $d.Content
    $sr.Find.ClearFormatting()
    $sr.Find.Execute($marker, $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
    $sr.Text = $replacement
    $sr.Bold = 1
    $sr.Bold = 0
}

# ---------------------------------------------------------------------------
# 1) Rewrite the first TECHNICAL SKILLS bullet.
#    Old:  "JS, HTML, CSS, React+Redux, GitHub, " + "Webpack, Sass"
#    New:  "JavaScript, HTML, CSS, React, Redux, Sass, GitHub, Webpack, BEM"
#          + ", " + "ООП" + ", " + "Redux-thunk"
# ---------------------------------------------------------------------------
$skillsPara = $d.Paragraphs.Item(30)
$skillsRange = $d.Range($skillsPara.Range.Start, $skillsPara.Range.End - 1)
$skillsRange.Text = "JavaScript, HTML, CSS, React, Redux, Sass, GitHub, Webpack, BEM#MARK1#OOP#MARK2#Redux-thunk"

Split-Marker "#MARK1#" ", "
Split-Marker "OOP" "ООП"
Split-Marker "#MARK2#" ", "

# ---------------------------------------------------------------------------
# 2) Insert two brand-new bulleted paragraphs right after it, reusing the
#    same list (numId 28) / paragraph formatting by letting
#    InsertParagraphAfter clone the source paragraph's pPr + rPr.
# ---------------------------------------------------------------------------
$skillsPara = $d.Paragraphs.Item(30)
$skillsPara.Range.InsertParagraphAfter()

$paraB = $d.Paragraphs.Item(31)
$paraBRange = $d.Range($paraB.Range.Start, $paraB.Range.End - 1)
$paraBRange.Text = "Axios#M1#React-redux#M2#Redux-form#M3#Reselect#M4#Classnames#M5#React-router-dom"

Split-Marker "#M1#" ", "
Split-Marker "#M2#" ", "
Split-Marker "#M3#" ", "
Split-Marker "#M4#" ", "
Split-Marker "#M5#" ", "

$paraB = $d.Paragraphs.Item(31)
$paraB.Range.InsertParagraphAfter()

$paraC = $d.Paragraphs.Item(32)
$paraCRange = $d.Range($paraC.Range.Start, $paraC.Range.End - 1)
$paraCRange.Text = "REST#N1#FLUX#N2#SPA"

Split-Marker "#N1#" ", "
Split-Marker "#N2#" ", "

# ---------------------------------------------------------------------------
# 3) The paragraph that used to hold a single " " run right after the
#    skills list becomes completely empty (the run is deleted).
# ---------------------------------------------------------------------------
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $pp = $d.Paragraphs.Item($i)
    if ($pp.Range.Text -eq " ") {
        $pp.Range.Text = ""
        break
    }
}

# ---------------------------------------------------------------------------
# 4) Merge the three runs that spelled out
#       "I took a course in React " + "with Kuzyuberdin" + "."
#    into a single run "I took a course in React with Kuzyuberdin."
# ---------------------------------------------------------------------------
$sr = $d.Content
$sr.Find.ClearFormatting()
$sr.Find.Execute("I took a course in React with Kuzyuberdin.", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
if (-not $sr.Find.Found) {
    $sr2 = $d.Content
    $sr2.Find.ClearFormatting()
    $sr2.Find.Execute("I took a course in React with Kuzyuberdin", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
}

Write-Output "done"
